# Daily attendance processing - 2025-10-07 18:26:00
# For every row's "Recorded By" (column G) value that contains multiple
# comma-separated recorders, rotate the list left by one position (the
# first recorder moves to the end of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rest = $parts[1..($parts.Count - 1)]
            $newParts = $rest + @($parts[0])
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value2 = $newVal
        }
    }
}
